# Add a new "3d_classic" worksheet right after the existing "3d" sheet.
# It is the same data as "3d" but with the 2nd axis header changed from
# "b\c" (the combined/new-style axis name) to "b" (the classic name),
# matching the existing "2d" / "2d_classic" pairing already in the
# workbook.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("3d")

# Copy "3d" and place the copy immediately after it.
$src.Copy($null, $src)

# Excel names the copy "3d (2)" by default; rename + tweak its header.
$newSheet = $wb.Worksheets.Item("3d (2)")
$newSheet.Name = "3d_classic"
$newSheet.Range("B1").Value = "b"

# Make the newly added sheet the active one (mirrors Excel's behaviour
# of selecting a sheet right after it is created/copied).
$newSheet.Activate()
